# Update cryptocurrency price and volume figures (Thu Nov  9 06:59:51 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '36.760.62'
Set-TextCell 2 5 '  +4.17%  '

Set-TextCell 3 4 '1.927.45'
Set-TextCell 3 5 '  +2.54%  '

Set-TextCell 4 5 '  -0.02%  '

Set-TextCell 5 4 '250.47'
Set-TextCell 5 5 '  +1.52%  '

Set-TextCell 6 4 '0.702'
Set-TextCell 6 5 '  +2.69%  '

Set-TextCell 7 5 '  +0.01%  '

Set-TextCell 8 4 '44.20'
Set-TextCell 8 5 '  +1.40%  '

Set-TextCell 9 4 '58.82'
Set-TextCell 9 5 '  +9.53%  '

Set-TextCell 10 5 '  +4.01%  '

Set-TextCell 11 5 '  +3.84%  '

Set-TextCell 12 5 '  +2.80%  '

Set-TextCell 13 4 '14.75'
Set-TextCell 13 5 '  +8.64%  '

Set-TextCell 14 5 '  +7.62%  '

Set-TextCell 15 4 '2.208.06'
Set-TextCell 15 5 '  +2.57%  '

Set-TextCell 16 5 '  +4.55%  '

Set-TextCell 17 4 '1.923.51'
Set-TextCell 17 5 '  +2.34%  '

Set-TextCell 18 4 '36.757.26'
Set-TextCell 18 5 '  +4.03%  '

Set-TextCell 19 4 '74.72'
Set-TextCell 19 5 '  +2.70%  '

Set-TextCell 20 5 '  +5.45%  '

Set-TextCell 21 4 '252.45'
Set-TextCell 21 5 '  +3.42%  '

Set-TextCell 22 4 '13.44'
Set-TextCell 22 5 '  +4.68%  '

Set-TextCell 23 4 '5.27'
Set-TextCell 23 5 '  +5.74%  '

Set-TextCell 24 4 '2.69'
Set-TextCell 24 5 '  +2.12%  '

Set-TextCell 25 4 '0.999'
Set-TextCell 25 5 '  -0.10%  '

Set-TextCell 26 4 '2.22'
Set-TextCell 26 5 '  +1.05%  '

Set-TextCell 27 4 '168.26'
Set-TextCell 27 5 '  +1.67%  '

Set-TextCell 28 4 '8.89'
Set-TextCell 28 5 '  +3.77%  '

Set-TextCell 29 5 '  +3.04%  '

Set-TextCell 30 5 '  +2.25%  '

Set-TextCell 31 5 '  +6.45%  '

Set-TextCell 32 4 '0.0621'
Set-TextCell 32 5 '  +4.65%  '

Set-TextCell 33 5 '  -3.34%  '

Set-TextCell 34 4 '4.38'
Set-TextCell 34 5 '  +5.36%  '

Set-TextCell 35 5 '  -0.01%  '

Set-TextCell 36 4 '0.0868'
Set-TextCell 36 5 '  +20.93%  '

Set-TextCell 37 4 '1.52'
Set-TextCell 37 5 '  -11.24%  '

Set-TextCell 38 4 '0.905'
Set-TextCell 38 5 '  +8.07%  '

Set-TextCell 39 4 '17.81'
Set-TextCell 39 5 '  +49.28%  '

Set-TextCell 40 4 '2.07'
Set-TextCell 40 5 '  +6.59%  '

Set-TextCell 41 4 '106.99'
Set-TextCell 41 5 '  +11.04%  '

Set-TextCell 42 5 '  +5.13%  '

Set-TextCell 43 4 '17.35'
Set-TextCell 43 5 '  -1.65%  '

Set-TextCell 44 5 '  +3.84%  '

Set-TextCell 45 4 '1.344.43'
Set-TextCell 45 5 '  +3.07%  '

Set-TextCell 46 4 '2.61'
Set-TextCell 46 5 '  +9.56%  '

Set-TextCell 47 5 '  +1.55%  '

Set-TextCell 48 4 '0.0816'
Set-TextCell 48 5 '  +1.66%  '

Set-TextCell 49 4 '2.80'
Set-TextCell 49 5 '  +2.81%  '

Set-TextCell 50 4 '6.47'
Set-TextCell 50 5 '  +3.87%  '

Set-TextCell 51 4 '43.55'
Set-TextCell 51 5 '  +3.53%  '

